# Add two new columns, I (I0) and J (IF), to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (bold font, border, centered alignment) from the
# existing header cell H1 onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data for rows 2..41: I (I0) and J (IF)
$data = @(
    @(8,8),
    @(4,6),
    @(7,8),
    @(7,7),
    @(7,7),
    @(7,8),
    @(4,5),
    @(7,7),
    @(8,8),
    @(7,7),
    @(10,10),
    @(7,8),
    @(9,9),
    @(8,8),
    @(7,8),
    @(8,8),
    @(8,8),
    @(8,9),
    @(5,6),
    @(7,8),
    @(7,8),
    @(7,9),
    @(6,8),
    @(6,7),
    @(7,7),
    @(5,8),
    @(8,9),
    @(5,6),
    @(8,9),
    @(1,5),
    @(1,5),
    @(1,6),
    @(1,6),
    @(1,5),
    @(1,6),
    @(1,6),
    @(1,4),
    @(4,6),
    @(8,8),
    @(1,2)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
